# Refresh the cryptos price/volume table with the latest scrape.
# (Row 44/45 also swap rank: InjectiveProtocol now outranks EnergySwap.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.375.59"
$ws.Range("E2").Value = "  +2.12%  "

# Row 3
$ws.Range("D3").Value = "2.674.64"
$ws.Range("E3").Value = "  +3.14%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'579.15"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6
$ws.Range("D6").Value = "'144.69"
$ws.Range("E6").Value = "  +1.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +0.72%  "

# Row 9
$ws.Range("E9").Value = "  +1.37%  "

# Row 10
$ws.Range("E10").Value = "  +4.48%  "

# Row 11
$ws.Range("D11").Value = "'0.382"
$ws.Range("E11").Value = "  +3.23%  "

# Row 12
$ws.Range("E12").Value = "  +1.27%  "

# Row 13
$ws.Range("D13").Value = "3.122.77"
$ws.Range("E13").Value = "  +2.25%  "

# Row 14
$ws.Range("D14").Value = "'26.18"
$ws.Range("E14").Value = "  +6.17%  "

# Row 15
$ws.Range("D15").Value = "61.300.99"
$ws.Range("E15").Value = "  +1.97%  "

# Row 16
$ws.Range("D16").Value = "'0.0000146"
$ws.Range("E16").Value = "  +3.69%  "

# Row 17
$ws.Range("D17").Value = "2.665.32"
$ws.Range("E17").Value = "  +2.44%  "

# Row 18
$ws.Range("D18").Value = "'11.73"
$ws.Range("E18").Value = "  +2.05%  "

# Row 19
$ws.Range("D19").Value = "'4.78"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
$ws.Range("D20").Value = "'354.90"
$ws.Range("E20").Value = "  +2.58%  "

# Row 21
$ws.Range("D21").Value = "'6.89"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("D23").Value = "'0.528"
$ws.Range("E23").Value = "  +0.54%  "

# Row 24
$ws.Range("D24").Value = "'64.54"
$ws.Range("E24").Value = "  +2.50%  "

# Row 25
$ws.Range("E25").Value = "  +3.26%  "

# Row 26
$ws.Range("D26").Value = "'8.53"
$ws.Range("E26").Value = "  +6.33%  "

# Row 27
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.50%  "

# Row 28
$ws.Range("E28").Value = "  +7.23%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0823"
$ws.Range("E29").Value = "  +4.59%  "

# Row 30
$ws.Range("E30").Value = "  +8.51%  "

# Row 31
$ws.Range("E31").Value = "  +2.74%  "

# Row 32
$ws.Range("E32").Value = "  -0.10%  "

# Row 33
$ws.Range("D33").Value = "'20.16"
$ws.Range("E33").Value = "  +3.87%  "

# Row 34
$ws.Range("D34").Value = "'1.13"
$ws.Range("E34").Value = "  +14.94%  "

# Row 35
$ws.Range("E35").Value = "  +10.08%  "

# Row 36
$ws.Range("E36").Value = "  +11.01%  "

# Row 37
$ws.Range("D37").Value = "'1.72"
$ws.Range("E37").Value = "  +6.38%  "

# Row 38
$ws.Range("D38").Value = "'0.968"
$ws.Range("E38").Value = "  +15.45%  "

# Row 39
$ws.Range("D39").Value = "'341.50"
$ws.Range("E39").Value = "  +10.76%  "

# Row 40
$ws.Range("D40").Value = "'4.16"
$ws.Range("E40").Value = "  +6.71%  "

# Row 41
$ws.Range("D41").Value = "'38.45"
$ws.Range("E41").Value = "  +1.14%  "

# Row 42
$ws.Range("D42").Value = "'5.39"
$ws.Range("E42").Value = "  +7.93%  "

# Row 43
$ws.Range("E43").Value = "  +6.63%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'21.21"
$ws.Range("E44").Value = "  +6.19%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'20.68"
$ws.Range("E45").Value = "  +5.29%  "

# Row 46
$ws.Range("E46").Value = "  +4.79%  "

# Row 47
$ws.Range("D47").Value = "'135.73"
$ws.Range("E47").Value = "  +0.35%  "

# Row 48
$ws.Range("D48").Value = "'0.0253"
$ws.Range("E48").Value = "  +5.01%  "

# Row 49
$ws.Range("E49").Value = "  +1.39%  "

# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").Value = "2.100.36"
$ws.Range("E51").Value = "  +4.19%  "
